$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 8, 4, 12),
    @(3, 7, 4, 13),
    @(3, 5, 5, 15),
    @(5, 7, 6, 13),
    @(5, 4, 6, 16),
    @(3, 7, 4, 13),
    @(5, 7, 3, 13),
    @(5, 5, 4, 15),
    @(5, 12, 6, 8),
    @(6, 6, 2, 14),
    @(2, 7, 3, 13),
    @(7, 5, 4, 15),
    @(4, 14, 3, 6)
)

$startRow = 1110
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("A1123").Select()
